$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for the rows whose data was repulled.
$ws.Range("F6").Value = 2
$ws.Range("F10").Value = -3
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = 0
$ws.Range("F16").Value = -2
